$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 ---------------------------------------------------------
$ws.Range("B2").Value = 312312
$ws.Range("C2").Value = "Абоба"
$ws.Range("D2").Value = "25.01.2024"
$ws.Range("E2").Value = "29.01.2024"

# --- Clear row 3 entirely (content + formatting) --------------------------
$ws.Range("A3:F3").Clear()

# --- Update row 4 ---------------------------------------------------------
$ws.Range("B4").Value = 123123123
$ws.Range("C4").Value = "ававпвпв"
$ws.Range("D4").Value = "25.01.2024"
$ws.Range("E4").Value = "26.01.2024"

# --- Clear row 5 entirely (content + formatting) --------------------------
$ws.Range("A5:F5").Clear()

# --- Drop the now-unused column F (days-until-overdue count) --------------
$ws.Range("F2").Clear()
$ws.Range("F4").Clear()

# --- Selection / cursor position updates as recorded in the sheet view ----
$ws.Range("H8").Select()
